$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("KPI")

# --- Row 27: section title (bold label, no border) ---
$c = $ws2.Range("B27")
$c.Font.Size = 9
$c.Font.Bold = $true
$c.Value = "Métricas que estão vendendo para esse produto"

# --- Column B labels, rows 28-38 (bold label, bordered) ---
$labels = @{
    28 = "CPV"
    29 = "% Cartões"
    30 = "%Boletos"
    31 = "Conversão Boletos"
    32 = "Taxa Conversão"
    33 = "CTR"
    34 = "Cliques/Venda"
    35 = "CPC"
    36 = "CPM"
    37 = "ROI"
    38 = "Investimento"
}
foreach ($r in 28..38) {
    $c = $ws2.Range("B$r")
    $c.Font.Size = 9
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.Value = $labels[$r]
}

# --- Column C percentage cells (rows 29-32), bordered, "0%" format ---
$pctVals = @{
    29 = "cartoes_venda"
    30 = "boletos_venda"
    31 = "conv_boleto_venda"
}
foreach ($r in 29,30,31) {
    $c = $ws2.Range("C$r")
    $c.NumberFormat = "0%"
    $c.Borders.LineStyle = 1
    $c.Value = $pctVals[$r]
}
# C32 gets the same percent style via a formula instead of a literal value
$c = $ws2.Range("C32")
$c.NumberFormat = "0%"
$c.Borders.LineStyle = 1

# --- Column C decimal cells (0.00 format), bordered ---
$decVals = @{
    28 = "cpv_venda"
    33 = "ctr_venda"
    35 = "cpc_venda"
    36 = "cpm_venda"
    37 = "roi_venda"
    38 = "spend_venda"
}
foreach ($r in 28,33,35,36,37,38) {
    $c = $ws2.Range("C$r")
    $c.NumberFormat = "0.00"
    $c.Borders.LineStyle = 1
    $c.Value = $decVals[$r]
}

# --- Column C integer cell (row 34), bordered, "0" format ---
$c = $ws2.Range("C34")
$c.NumberFormat = "0"
$c.Borders.LineStyle = 1
$c.Value = "clpv_venda"

# --- Formula for C32 (taxa de conversão) ---
$ws2.Range("C32").Formula = "=IF(((1/C34)*C29)+((1/C34)*C30*C31)= 0, 0.005, ((1/C34)*C29)+((1/C34)*C30*C31))"

# --- Selection / activation: KPI sheet becomes the active tab, C34 selected ---
$ws2.Range("C34").Select() | Out-Null
$ws2.Activate()
